# Update the crypto price/volume table on Sheet1 to match the latest
# GitHub Actions scrape (commit: "Updated cryptos list on Sun Mar 3
# 05:57:31 UTC 2024 with GitHub Actions").
#
# Price values in column D are stored as text (e.g. "61.999.03",
# "43.60") so that Excel does not reinterpret them as numbers/dates.
# For cells whose new text looks like a plain number we prefix the
# value with a leading apostrophe to force text entry, then reset the
# cell style to "Normal" so no stray number-format styling is left
# behind once Excel applies a text format for the apostrophe-prefixed
# entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '''61.952.92'
$ws.Cells.Item(2,4).Style = 'Normal'
$ws.Cells.Item(2,5).Value = '  -0.06%  '
$ws.Cells.Item(3,4).Value = '''3.423.89'
$ws.Cells.Item(3,4).Style = 'Normal'
$ws.Cells.Item(3,5).Value = '  +0.11%  '
$ws.Cells.Item(4,5).Value = '  +0.03%  '
$ws.Cells.Item(5,4).Value = '''411.16'
$ws.Cells.Item(5,4).Style = 'Normal'
$ws.Cells.Item(5,5).Value = '  +0.43%  '
$ws.Cells.Item(6,4).Value = '''129.70'
$ws.Cells.Item(6,4).Style = 'Normal'
$ws.Cells.Item(6,5).Value = '  +1.03%  '
$ws.Cells.Item(7,5).Value = '  +1.44%  '
$ws.Cells.Item(8,5).Value = '  +0.06%  '
$ws.Cells.Item(9,4).Value = '''0.727'
$ws.Cells.Item(9,4).Style = 'Normal'
$ws.Cells.Item(9,5).Value = '  -2.75%  '
$ws.Cells.Item(10,4).Value = '''0.140'
$ws.Cells.Item(10,4).Style = 'Normal'
$ws.Cells.Item(10,5).Value = '  -0.94%  '
$ws.Cells.Item(11,4).Value = '''43.55'
$ws.Cells.Item(11,4).Style = 'Normal'
$ws.Cells.Item(11,5).Value = '  +1.44%  '
$ws.Cells.Item(12,4).Value = '''9.23'
$ws.Cells.Item(12,4).Style = 'Normal'
$ws.Cells.Item(12,5).Value = '  +4.51%  '
$ws.Cells.Item(13,4).Value = '''0.0000216'
$ws.Cells.Item(13,4).Style = 'Normal'
$ws.Cells.Item(13,5).Value = '  +7.31%  '
$ws.Cells.Item(14,4).Value = '''3.970.30'
$ws.Cells.Item(14,4).Style = 'Normal'
$ws.Cells.Item(14,5).Value = '  +0.27%  '
$ws.Cells.Item(16,4).Value = '''21.09'
$ws.Cells.Item(16,4).Style = 'Normal'
$ws.Cells.Item(16,5).Value = '  -0.16%  '
$ws.Cells.Item(17,4).Value = '''3.427.83'
$ws.Cells.Item(17,4).Style = 'Normal'
$ws.Cells.Item(17,5).Value = '  +0.78%  '
$ws.Cells.Item(18,4).Value = '''1.10'
$ws.Cells.Item(18,4).Style = 'Normal'
$ws.Cells.Item(18,5).Value = '  +4.06%  '
$ws.Cells.Item(19,4).Value = '''12.51'
$ws.Cells.Item(19,4).Style = 'Normal'
$ws.Cells.Item(19,5).Value = '  -1.55%  '
$ws.Cells.Item(20,4).Value = '''61.939.15'
$ws.Cells.Item(20,4).Style = 'Normal'
$ws.Cells.Item(20,5).Value = '  -0.17%  '
$ws.Cells.Item(21,4).Value = '''491.51'
$ws.Cells.Item(21,4).Style = 'Normal'
$ws.Cells.Item(21,5).Value = '  +21.96%  '
$ws.Cells.Item(22,4).Value = '''92.46'
$ws.Cells.Item(22,4).Style = 'Normal'
$ws.Cells.Item(22,5).Value = '  +2.45%  '
$ws.Cells.Item(23,5).Value = '  +4.30%  '
$ws.Cells.Item(24,4).Value = '''13.68'
$ws.Cells.Item(24,4).Style = 'Normal'
$ws.Cells.Item(24,5).Value = '  +1.60%  '
$ws.Cells.Item(25,4).Value = '''3.37'
$ws.Cells.Item(25,4).Style = 'Normal'
$ws.Cells.Item(25,5).Value = '  +4.28%  '
$ws.Cells.Item(26,4).Value = '''34.56'
$ws.Cells.Item(26,4).Style = 'Normal'
$ws.Cells.Item(26,5).Value = '  +5.11%  '
$ws.Cells.Item(27,4).Value = '''9.13'
$ws.Cells.Item(27,4).Style = 'Normal'
$ws.Cells.Item(27,5).Value = '  +6.80%  '
$ws.Cells.Item(28,4).Value = '''4.81'
$ws.Cells.Item(28,4).Style = 'Normal'
$ws.Cells.Item(28,5).Value = '  +0.38%  '
$ws.Cells.Item(29,4).Value = '''7.79'
$ws.Cells.Item(29,4).Style = 'Normal'
$ws.Cells.Item(29,5).Value = '  +1.83%  '
$ws.Cells.Item(30,5).Value = '  -0.05%  '
$ws.Cells.Item(31,4).Value = '''12.07'
$ws.Cells.Item(31,4).Style = 'Normal'
$ws.Cells.Item(31,5).Value = '  +2.12%  '
$ws.Cells.Item(32,2).Value = 'Hedera'
$ws.Cells.Item(32,3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(32,4).Value = '''0.115'
$ws.Cells.Item(32,4).Style = 'Normal'
$ws.Cells.Item(32,5).Value = '  -3.16%  '
$ws.Cells.Item(33,2).Value = 'Kaspa'
$ws.Cells.Item(33,3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(33,4).Value = '''0.168'
$ws.Cells.Item(33,4).Style = 'Normal'
$ws.Cells.Item(33,5).Value = '  -2.03%  '
$ws.Cells.Item(34,4).Value = '''42.09'
$ws.Cells.Item(34,4).Style = 'Normal'
$ws.Cells.Item(34,5).Value = '  -4.22%  '
$ws.Cells.Item(35,5).Value = '  -0.01%  '
$ws.Cells.Item(36,4).Value = '''58.12'
$ws.Cells.Item(36,4).Style = 'Normal'
$ws.Cells.Item(36,5).Value = '  +9.84%  '
$ws.Cells.Item(37,4).Value = '''0.0496'
$ws.Cells.Item(37,4).Style = 'Normal'
$ws.Cells.Item(37,5).Value = '  -0.28%  '
$ws.Cells.Item(38,5).Value = '  -0.02%  '
$ws.Cells.Item(39,4).Value = '''3.45'
$ws.Cells.Item(39,4).Style = 'Normal'
$ws.Cells.Item(39,5).Value = '  +1.66%  '
$ws.Cells.Item(40,2).Value = 'TheGraph'
$ws.Cells.Item(40,3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(40,4).Value = '''0.330'
$ws.Cells.Item(40,4).Style = 'Normal'
$ws.Cells.Item(40,5).Value = '  +4.57%  '
$ws.Cells.Item(41,2).Value = 'Monero'
$ws.Cells.Item(41,3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(41,4).Value = '''150.95'
$ws.Cells.Item(41,4).Style = 'Normal'
$ws.Cells.Item(41,5).Value = '  +7.06%  '
$ws.Cells.Item(42,5).Value = '  +3.64%  '
$ws.Cells.Item(43,2).Value = 'Stacks'
$ws.Cells.Item(43,3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(43,4).Value = '''2.98'
$ws.Cells.Item(43,4).Style = 'Normal'
$ws.Cells.Item(43,5).Value = '  +2.60%  '
$ws.Cells.Item(44,2).Value = 'ARBITRUM'
$ws.Cells.Item(44,3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(44,4).Value = '''2.14'
$ws.Cells.Item(44,4).Style = 'Normal'
$ws.Cells.Item(44,5).Value = '  +8.56%  '
$ws.Cells.Item(45,5).Value = '  +14.31%  '
$ws.Cells.Item(47,2).Value = 'Celestia'
$ws.Cells.Item(47,3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Cells.Item(47,4).Value = '''16.66'
$ws.Cells.Item(47,4).Style = 'Normal'
$ws.Cells.Item(47,5).Value = '  -0.68%  '
$ws.Cells.Item(48,2).Value = 'ThetaToken'
$ws.Cells.Item(48,3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(48,4).Value = '''2.34'
$ws.Cells.Item(48,4).Style = 'Normal'
$ws.Cells.Item(48,5).Value = '  +20.52%  '
$ws.Cells.Item(49,4).Value = '''22.82'
$ws.Cells.Item(49,4).Style = 'Normal'
$ws.Cells.Item(49,5).Value = '  +5.11%  '
$ws.Cells.Item(50,4).Value = '''117.77'
$ws.Cells.Item(50,4).Style = 'Normal'
$ws.Cells.Item(50,5).Value = '  +23.48%  '
$ws.Cells.Item(51,4).Value = '''0.146'
$ws.Cells.Item(51,4).Style = 'Normal'
$ws.Cells.Item(51,5).Value = '  +14.17%  '
